$d = $word.ActiveDocument

# --- Change 1: remove the redundant "the " in "a projection of the this year's..." ---
$d.Content.Find.Execute(
    "a projection of the this year" + [char]8217 + "s annual revenue",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "a projection of this year" + [char]8217 + "s annual revenue",
    2)

# --- Change 2: fix typo "would't" -> "wouldn't" by inserting the missing "n" ---
$find = $d.Content.Duplicate
$find.Find.Execute(
    "therefore it would" + [char]8217 + "t make sense",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Position right after "would" (before the apostrophe) where the "n" is missing.
$insertPos = $find.Start + "therefore it would".Length
$ins = $d.Range($insertPos, $insertPos)
$ins.InsertAfter("n")

# Toggling formatting on the newly inserted run forces Word to keep it as its
# own run (matching the author's original edit) instead of silently
# re-merging it with the neighboring text that shares identical formatting.
$insRun = $d.Range($insertPos, $insertPos + 1)
$insRun.Font.Bold = $true
$insRun.Font.Bold = $false
